$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsRacp  = $wb.Worksheets.Item("RACP")

# Update the calibration factor and its label on the About sheet.
$wsAbout.Range("A11").Value = 0.75350342301658668
$wsAbout.Range("B11").Value = "2023 dollars per 2012 dollar"

# Replace the formula in RACP!B2 with a hard-coded value.
$wsRacp.Range("B2").Value = 999

# Update view state: About sheet selection, then make RACP the active/selected sheet.
$wsAbout.Activate()
$wsAbout.Range("B12").Select()

$wsRacp.Activate()
$wsRacp.Range("E8").Select()
